$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2444105972001947
$ws.Range("C2").Value = 0.4739145320991289
$ws.Range("D2").Value = 0.3556655416862647
$ws.Range("E2").Value = 0.5963770130431459
$ws.Range("F2").Value = 0.5645288044211415
$ws.Range("G2").Value = 14
